# Productos.xlsx - add "referencia_producto"/"referencia_molde" reference columns
# Original columns: referencia | producto | molde | cantidad
# Target columns:   referencia_producto | producto | referencia_molde | molde | cantidad

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Make room: insert a new blank column before "molde" (current column C).
#    This shifts molde -> D and cantidad -> E.
# ---------------------------------------------------------------------------
$ws.Range("C1").EntireColumn.Insert()

# 2. Grow the table definition so it covers the new column too.
$lo.Resize($ws.Range("A1:E6"))

# ---------------------------------------------------------------------------
# 3. Rename / fill in the header row (also renames the table columns and
#    updates the shared-string table automatically).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "referencia_producto"
$ws.Range("B1").Value = "producto"
$ws.Range("C1").Value = "referencia_molde"
$ws.Range("D1").Value = "molde"
$ws.Range("E1").Value = "cantidad"

# ---------------------------------------------------------------------------
# Helper functions for thin black borders (matching the workbook's existing
# "boxed" table look).
# ---------------------------------------------------------------------------
function Set-ThinEdge($rng, $edge) {
    $rng.Borders.Item($edge).LineStyle = 1
    $rng.Borders.Item($edge).ColorIndex = 1
}
function Clear-Edge($rng, $edge) {
    $rng.Borders.Item($edge).LineStyle = -4142
}

$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlCenter = -4108

# ---------------------------------------------------------------------------
# 4. Style the new "referencia_molde" column (C) exactly like the "molde"
#    header/data look the workbook already uses (bold Arial 14 header with a
#    left/right/bottom border, Arial 12 boxed data cells).
# ---------------------------------------------------------------------------
$c1 = $ws.Range("C1")
$c1.Font.Name = "Arial"
$c1.Font.Size = 14
$c1.Font.Bold = $true
$c1.HorizontalAlignment = $xlCenter
$c1.VerticalAlignment = $xlCenter
Set-ThinEdge $c1 $xlEdgeLeft
Set-ThinEdge $c1 $xlEdgeRight
Set-ThinEdge $c1 $xlEdgeBottom
Clear-Edge $c1 $xlEdgeTop

$c2 = $ws.Range("C2")
$c2.Font.Name = "Arial"
$c2.Font.Size = 12
$c2.HorizontalAlignment = $xlCenter
$c2.VerticalAlignment = $xlCenter
Set-ThinEdge $c2 $xlEdgeLeft
Set-ThinEdge $c2 $xlEdgeRight
Set-ThinEdge $c2 $xlEdgeTop
Clear-Edge $c2 $xlEdgeBottom

$cMid = $ws.Range("C3:C5")
$cMid.Font.Name = "Arial"
$cMid.Font.Size = 12
$cMid.HorizontalAlignment = $xlCenter
$cMid.VerticalAlignment = $xlCenter
Set-ThinEdge $cMid $xlEdgeLeft
Set-ThinEdge $cMid $xlEdgeRight
Clear-Edge $cMid $xlEdgeTop
Clear-Edge $cMid $xlEdgeBottom

$c6 = $ws.Range("C6")
$c6.Font.Name = "Arial"
$c6.Font.Size = 12
$c6.HorizontalAlignment = $xlCenter
$c6.VerticalAlignment = $xlCenter
Set-ThinEdge $c6 $xlEdgeLeft
Set-ThinEdge $c6 $xlEdgeRight
Set-ThinEdge $c6 $xlEdgeBottom
Clear-Edge $c6 $xlEdgeTop

# ---------------------------------------------------------------------------
# 5. Style the new "cantidad" column (E, was D) with the same boxed look.
# ---------------------------------------------------------------------------
$e1 = $ws.Range("E1")
$e1.Font.Name = "Arial"
$e1.Font.Size = 14
$e1.Font.Bold = $true
$e1.HorizontalAlignment = $xlCenter
$e1.VerticalAlignment = $xlCenter
Set-ThinEdge $e1 $xlEdgeLeft
Set-ThinEdge $e1 $xlEdgeRight
Set-ThinEdge $e1 $xlEdgeBottom
Clear-Edge $e1 $xlEdgeTop

$e2 = $ws.Range("E2")
$e2.Font.Name = "Arial"
$e2.Font.Size = 12
$e2.HorizontalAlignment = $xlCenter
$e2.VerticalAlignment = $xlCenter
Set-ThinEdge $e2 $xlEdgeLeft
Set-ThinEdge $e2 $xlEdgeRight
Set-ThinEdge $e2 $xlEdgeTop
Clear-Edge $e2 $xlEdgeBottom

$eMid = $ws.Range("E3:E5")
$eMid.Font.Name = "Arial"
$eMid.Font.Size = 12
$eMid.HorizontalAlignment = $xlCenter
$eMid.VerticalAlignment = $xlCenter
Set-ThinEdge $eMid $xlEdgeLeft
Set-ThinEdge $eMid $xlEdgeRight
Clear-Edge $eMid $xlEdgeTop
Clear-Edge $eMid $xlEdgeBottom

$e6 = $ws.Range("E6")
$e6.Font.Name = "Arial"
$e6.Font.Size = 12
$e6.HorizontalAlignment = $xlCenter
$e6.VerticalAlignment = $xlCenter
Set-ThinEdge $e6 $xlEdgeLeft
Set-ThinEdge $e6 $xlEdgeRight
Set-ThinEdge $e6 $xlEdgeBottom
Clear-Edge $e6 $xlEdgeTop

# ---------------------------------------------------------------------------
# 6. Column widths / dimension housekeeping.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 28.57
$ws.Columns.Item(5).ColumnWidth = 19.71

# ---------------------------------------------------------------------------
# 7. Selection, matching what the author left selected after the edit.
# ---------------------------------------------------------------------------
$ws.Range("D12").Select()
